# Update cryptocurrency price/volume data (automated refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "26.658.66"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "1.633.89"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'212.54"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "'0.493"
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("D10").Value = "'19.03"
$ws.Range("E10").Value = "  +3.39%  "
$ws.Range("D11").Value = "'0.0840"
$ws.Range("E11").Value = "  +3.37%  "
$ws.Range("D12").Value = "1.862.40"
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("D13").Value = "1.653.85"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("D16").Value = "26.663.12"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("D17").Value = "'62.93"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("D19").Value = "'208.97"
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").Value = "'6.17"
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "'1.92"
$ws.Range("E24").Value = "  +2.18%  "
$ws.Range("D25").Value = "'146.54"
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("D28").Value = "'6.76"
$ws.Range("E28").Value = "  +2.86%  "
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "1.168.59"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").Value = "'0.0167"
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("E38").Value = "  +2.83%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'0.502"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D44").Value = "1.774.57"
$ws.Range("E44").Value = "  +2.00%  "
$ws.Range("D45").Value = "'91.93"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("E47").Value = "  -1.25%  "
$ws.Range("D48").Value = "'54.62"
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("E49").Value = "  +1.43%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.409"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.51"
$ws.Range("E51").Value = "  +4.02%  "
